$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp header (refreshed snapshot time)
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 28 de Marzo de 2020 a las 22:29'

# Refreshed COVID-19 country case counts. Some countries crossed rank
# boundaries versus their neighbours once totals were updated, so the
# country label in col A for a handful of rows changes along with the data.

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = 'Estados Unidos'
$ws.Cells.Item(4, 2).Value = 120204
$ws.Cells.Item(4, 3).Value = 16078
$ws.Cells.Item(4, 4).Value = 3229
$ws.Cells.Item(4, 5).Value = 114978
$ws.Cells.Item(4, 6).Value = 2666
$ws.Cells.Item(4, 7).Value = 301
$ws.Cells.Item(4, 8).Value = 1997

# Row 8: Alemania
$ws.Cells.Item(8, 1).Value = 'Alemania'
$ws.Cells.Item(8, 2).Value = 57695
$ws.Cells.Item(8, 3).Value = 6824
$ws.Cells.Item(8, 4).Value = 8481
$ws.Cells.Item(8, 5).Value = 48781
$ws.Cells.Item(8, 6).Value = 1581
$ws.Cells.Item(8, 7).Value = 82
$ws.Cells.Item(8, 8).Value = 433

# Row 18: Canada
$ws.Cells.Item(18, 1).Value = 'Canada'
$ws.Cells.Item(18, 2).Value = 5576
$ws.Cells.Item(18, 3).Value = 819
$ws.Cells.Item(18, 4).Value = 354
$ws.Cells.Item(18, 5).Value = 5167
$ws.Cells.Item(18, 6).Value = 120
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 55

# Row 20: Noruega
$ws.Cells.Item(20, 1).Value = 'Noruega'
$ws.Cells.Item(20, 2).Value = 4012
$ws.Cells.Item(20, 3).Value = 241
$ws.Cells.Item(20, 4).Value = 7
$ws.Cells.Item(20, 5).Value = 3982
$ws.Cells.Item(20, 6).Value = 84
$ws.Cells.Item(20, 7).Value = 4
$ws.Cells.Item(20, 8).Value = 23

# Row 35: Rumania
$ws.Cells.Item(35, 1).Value = 'Rumania'
$ws.Cells.Item(35, 2).Value = 1452
$ws.Cells.Item(35, 3).Value = 160
$ws.Cells.Item(35, 4).Value = 139
$ws.Cells.Item(35, 5).Value = 1279
$ws.Cells.Item(35, 6).Value = 34
$ws.Cells.Item(35, 7).Value = 8
$ws.Cells.Item(35, 8).Value = 34

# Row 44: India
$ws.Cells.Item(44, 1).Value = 'India'
$ws.Cells.Item(44, 2).Value = 987
$ws.Cells.Item(44, 3).Value = 100
$ws.Cells.Item(44, 4).Value = 84
$ws.Cells.Item(44, 5).Value = 879
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 4
$ws.Cells.Item(44, 8).Value = 24

# Row 45: Islandia
$ws.Cells.Item(45, 1).Value = 'Islandia'
$ws.Cells.Item(45, 2).Value = 963
$ws.Cells.Item(45, 3).Value = 73
$ws.Cells.Item(45, 4).Value = 114
$ws.Cells.Item(45, 5).Value = 847
$ws.Cells.Item(45, 6).Value = 18
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 2

# Row 69: Marruecos
$ws.Cells.Item(69, 1).Value = 'Marruecos'
$ws.Cells.Item(69, 2).Value = 390
$ws.Cells.Item(69, 3).Value = 45
$ws.Cells.Item(69, 4).Value = 11
$ws.Cells.Item(69, 5).Value = 354
$ws.Cells.Item(69, 6).Value = 1
$ws.Cells.Item(69, 7).Value = 2
$ws.Cells.Item(69, 8).Value = 25

# Row 70: Ucrania
$ws.Cells.Item(70, 1).Value = 'Ucrania'
$ws.Cells.Item(70, 2).Value = 356
$ws.Cells.Item(70, 3).Value = 46
$ws.Cells.Item(70, 4).Value = 5
$ws.Cells.Item(70, 5).Value = 342
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 4
$ws.Cells.Item(70, 8).Value = 9

# Row 71: Hungria
$ws.Cells.Item(71, 1).Value = 'Hungria'
$ws.Cells.Item(71, 2).Value = 343
$ws.Cells.Item(71, 3).Value = 43
$ws.Cells.Item(71, 4).Value = 34
$ws.Cells.Item(71, 5).Value = 298
$ws.Cells.Item(71, 6).Value = 6
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = 11

# Row 72: Bulgaria
$ws.Cells.Item(72, 1).Value = 'Bulgaria'
$ws.Cells.Item(72, 2).Value = 331
$ws.Cells.Item(72, 3).Value = 38
$ws.Cells.Item(72, 4).Value = 11
$ws.Cells.Item(72, 5).Value = 314
$ws.Cells.Item(72, 6).Value = 8
$ws.Cells.Item(72, 7).Value = 3
$ws.Cells.Item(72, 8).Value = 6

# Row 166: Islas Caimanes
$ws.Cells.Item(166, 1).Value = 'Islas Caimanes'
$ws.Cells.Item(166, 2).Value = 8
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 4).Value = 0
$ws.Cells.Item(166, 5).Value = 7
$ws.Cells.Item(166, 6).Value = 0
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 1

# Row 167: Guyana
$ws.Cells.Item(167, 1).Value = 'Guyana'
$ws.Cells.Item(167, 2).Value = 8
$ws.Cells.Item(167, 3).Value = 3
$ws.Cells.Item(167, 4).Value = 0
$ws.Cells.Item(167, 5).Value = 7
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 1

# Row 170: Antigua y Barbuda
$ws.Cells.Item(170, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(170, 2).Value = 7
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(170, 4).Value = 0
$ws.Cells.Item(170, 5).Value = 7
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 0
$ws.Cells.Item(170, 8).Value = 0

# Row 171: Granada
$ws.Cells.Item(171, 1).Value = 'Granada'
$ws.Cells.Item(171, 2).Value = 7
$ws.Cells.Item(171, 3).Value = 0
$ws.Cells.Item(171, 4).Value = 0
$ws.Cells.Item(171, 5).Value = 7
$ws.Cells.Item(171, 6).Value = 0
$ws.Cells.Item(171, 7).Value = 0
$ws.Cells.Item(171, 8).Value = 0

# Row 177: Mauritania
$ws.Cells.Item(177, 1).Value = 'Mauritania'
$ws.Cells.Item(177, 2).Value = 5
$ws.Cells.Item(177, 3).Value = 2
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 5).Value = 5
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 0

# Row 178: San Bartolome
$ws.Cells.Item(178, 1).Value = 'San Bartolome'
$ws.Cells.Item(178, 2).Value = 5
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(178, 4).Value = 0
$ws.Cells.Item(178, 5).Value = 5
$ws.Cells.Item(178, 6).Value = 0
$ws.Cells.Item(178, 7).Value = 0
$ws.Cells.Item(178, 8).Value = 0

# Row 179: Fiyi
$ws.Cells.Item(179, 1).Value = 'Fiyi'
$ws.Cells.Item(179, 2).Value = 5
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 0
$ws.Cells.Item(179, 5).Value = 5
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 0
$ws.Cells.Item(179, 8).Value = 0

# Row 180: Siria
$ws.Cells.Item(180, 1).Value = 'Siria'
$ws.Cells.Item(180, 2).Value = 5
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 5).Value = 5
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0

# Row 181: Montserrat
$ws.Cells.Item(181, 1).Value = 'Montserrat'
$ws.Cells.Item(181, 2).Value = 5
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 5
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 0

# Row 182: Angola
$ws.Cells.Item(182, 1).Value = 'Angola'
$ws.Cells.Item(182, 2).Value = 5
$ws.Cells.Item(182, 3).Value = 1
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 5
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 0

# Row 185: Sudan
$ws.Cells.Item(185, 1).Value = 'Sudan'
$ws.Cells.Item(185, 2).Value = 5
$ws.Cells.Item(185, 3).Value = 2
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 5).Value = 4
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 1

# Row 186: Islas Turcas y Caicos
$ws.Cells.Item(186, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(186, 2).Value = 4
$ws.Cells.Item(186, 3).Value = 2
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 4
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 0

# Row 187: Congo
$ws.Cells.Item(187, 1).Value = 'Congo'
$ws.Cells.Item(187, 2).Value = 4
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 5).Value = 4
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0
